$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$row = 54

# Copy the date-formatted style from the row above (style index 1) so we don't
# introduce brand new number-format styles for the new date cells.
$ws.Cells.Item($row - 1, 3).Copy() | Out-Null
$ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($row - 1, 4).Copy() | Out-Null
$ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = "What It Takes"
$ws.Cells.Item($row, 2).Value = "Stephen Schwarzman"
$ws.Cells.Item($row, 3).Value = 43931
$ws.Cells.Item($row, 4).Value = 43936
$ws.Cells.Item($row, 5).Value = "entreuprenuer;business;private equity;memoir;blackstone;investing;excellence"
$ws.Cells.Item($row, 6).Value = "Hard Copy"
$ws.Cells.Item($row, 7).Value = "354 Pages"

$ws.Range("A" + ($row + 1)).Select()
